$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.162.47'
$ws.Range('E2').Value = '  +3.07%  '
$ws.Range('D3').Value = '1.578.27'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.89'
$ws.Range('E5').Value = '  +1.14%  '
$ws.Range('E6').Value = '  +6.74%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.995'
$ws.Range('E7').Value = '  -0.59%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '25.91'
$ws.Range('E8').Value = '  +9.66%  '
$ws.Range('E9').Value = '  +2.53%  '
$ws.Range('E10').Value = '  +1.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0900'
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('D12').Value = '1.804.32'
$ws.Range('E12').Value = '  +1.79%  '
$ws.Range('D13').Value = '1.594.61'
$ws.Range('E13').Value = '  +2.65%  '
$ws.Range('D14').Value = '29.158.12'
$ws.Range('E14').Value = '  +3.04%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.70'
$ws.Range('E15').Value = '  +2.15%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.521'
$ws.Range('E16').Value = '  +2.45%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.34'
$ws.Range('E17').Value = '  +3.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '238.30'
$ws.Range('E18').Value = '  +4.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.44'
$ws.Range('E19').Value = '  +1.70%  '
$ws.Range('D20').Value = '0.0₃0692'
$ws.Range('E20').Value = '  +2.76%  '
$ws.Range('E21').Value = '  -0.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.99'
$ws.Range('E22').Value = '  +1.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.20'
$ws.Range('E23').Value = '  +4.21%  '
$ws.Range('E24').Value = '  +4.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.31'
$ws.Range('E25').Value = '  +2.74%  '
$ws.Range('E26').Value = '  +4.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.11'
$ws.Range('E27').Value = '  +2.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.32'
$ws.Range('E28').Value = '  +1.45%  '
$ws.Range('E29').Value = '  -0.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0465'
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('E31').Value = '  +0.58%  '
$ws.Range('E32').Value = '  +1.47%  '
$ws.Range('D33').Value = '1.418.71'
$ws.Range('E33').Value = '  +2.68%  '
$ws.Range('E34').Value = '  +0.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.05'
$ws.Range('E35').Value = '  -1.49%  '
$ws.Range('E36').Value = '  +1.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.75'
$ws.Range('E37').Value = '  +6.46%  '
$ws.Range('E38').Value = '  -2.13%  '
$ws.Range('E39').Value = '  +1.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.526'
$ws.Range('E40').Value = '  +3.37%  '
$ws.Range('B41').Value = 'BitcoinSV'
$ws.Range('C41').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '54.21'
$ws.Range('E41').Value = '  +28.57%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.95'
$ws.Range('E42').Value = '  +2.32%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.996'
$ws.Range('E43').Value = '  -0.55%  '
$ws.Range('E44').Value = '  +1.72%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0471'
$ws.Range('E45').Value = '  +1.30%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.55'
$ws.Range('E46').Value = '  +4.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.34'
$ws.Range('E47').Value = '  -0.69%  '
$ws.Range('D48').Value = '1.716.01'
$ws.Range('E48').Value = '  +1.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.849'
$ws.Range('E49').Value = '  -6.40%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '85.72'
$ws.Range('E50').Value = '  +0.47%  '
$ws.Range('E51').Value = '  +0.93%  '
